$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("L17").Value = 10880.7
$ws.Range("J17").Value = 3626.9
$ws.Range("N17").Value = -11216.7
$ws.Range("H17").Value = 5358.952
$ws.Range("L70").Value = 18675
$ws.Range("J70").Value = 6225
$ws.Range("N70").Value = -19215
$ws.Range("H70").Value = 3001126.2
$ws.Range("H73").Value = 3001126.2
$ws.Range("L73").Value = 18675
$ws.Range("N73").Value = -20547
$ws.Range("J73").Value = 6225
$ws.Range("K98").Value = 5333.778
$ws.Range("I98").Value = 5333.778
$ws.Range("N98").Value = -7848
$ws.Range("J98").Value = 4852
$ws.Range("M98").Value = -3835.778
$ws.Range("L98").Value = 4852
$ws.Range("H98").Value = 5213.3335
$ws.Range("L122").Value = 14556
$ws.Range("N122").Value = -19456
$ws.Range("I122").Value = 5333.778
$ws.Range("H122").Value = 5213.3335
$ws.Range("J122").Value = 4852
$ws.Range("M122").Value = -13551.334
$ws.Range("K122").Value = 16001.334
$ws.Range("M125").Value = -6082.5003
$ws.Range("H125").Value = 942.1429000000001
$ws.Range("K125").Value = 8542.5003
$ws.Range("I125").Value = 949.1667
$ws.Range("M137").Value = -2142
$ws.Range("J137").Value = 3500
$ws.Range("N137").Value = -15600
$ws.Range("K137").Value = 4692
$ws.Range("L137").Value = 10500
$ws.Range("H137").Value = 2532
$ws.Range("I137").Value = 1564
$ws.Range("L138").Value = 7211.7276
$ws.Range("K138").Value = 3388.8
$ws.Range("M138").Value = 1751.2
$ws.Range("J138").Value = 2403.9092
$ws.Range("H138").Value = 1518.9722
$ws.Range("N138").Value = -17491.7276
$ws.Range("I138").Value = 1129.6

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("M5").Value = -42.80000000000001
$ws.Range("I5").Value = 154.8
$ws.Range("K5").Value = 154.8
$ws.Range("H5").Value = 2795.2
$ws.Range("H32").Value = 2586.1875
$ws.Range("J32").Value = 1050.6666
$ws.Range("M32").Value = -2401.5557
$ws.Range("N32").Value = -1624.6666
$ws.Range("L32").Value = 1050.6666
$ws.Range("K32").Value = 2688.5557
$ws.Range("I32").Value = 2688.5557
$ws.Range("M61").Value = -20836020
$ws.Range("H61").Value = 11908187
$ws.Range("K61").Value = 20836232
$ws.Range("I61").Value = 20836232
$ws.Range("N63").Value = -28580158
$ws.Range("L63").Value = 28578786
$ws.Range("J63").Value = 28578786
$ws.Range("M63").Value = -500000314
$ws.Range("H63").Value = 133339280
$ws.Range("K63").Value = 500001000
$ws.Range("I63").Value = 500001000
$ws.Range("N66").Value = -142900794
$ws.Range("M66").Value = -2500001568
$ws.Range("I66").Value = 500001000
$ws.Range("L66").Value = 142893930
$ws.Range("J66").Value = 28578786
$ws.Range("K66").Value = 2500005000
$ws.Range("H66").Value = 133339280
$ws.Range("L122").Value = 11549.25
$ws.Range("N122").Value = -16449.25
$ws.Range("I122").Value = 17546576
$ws.Range("H122").Value = 12348731
$ws.Range("J122").Value = 3849.75
$ws.Range("M122").Value = -52637278
$ws.Range("K122").Value = 52639728
$ws.Range("H125").Value = 45000
$ws.Range("N125").Value = -54840
$ws.Range("J125").Value = 45000
$ws.Range("L125").Value = 45000
$ws.Range("K136").Value = 62508696
$ws.Range("I136").Value = 20836232
$ws.Range("H136").Value = 11908187
$ws.Range("M136").Value = -62506146

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("I4").Value = 154.8
$ws.Range("M4").Value = -39.80000000000001
$ws.Range("H4").Value = 2795.2
$ws.Range("K4").Value = 154.8
$ws.Range("K96").Value = 29975.666
$ws.Range("M96").Value = -27229.666
$ws.Range("I96").Value = 29975.666
$ws.Range("H96").Value = 29975.666

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("K31").Value = 1492.081
$ws.Range("I31").Value = 1492.081
$ws.Range("H31").Value = 2599.057
$ws.Range("M31").Value = -1197.081
$ws.Range("H34").Value = 2599.057
$ws.Range("I34").Value = 1492.081
$ws.Range("M34").Value = -1290.081
$ws.Range("K34").Value = 1492.081
$ws.Range("J62").Value = 333340670
$ws.Range("K62").Value = 6688.222
$ws.Range("H62").Value = 83340184
$ws.Range("M62").Value = -6064.222
$ws.Range("I62").Value = 6688.222
$ws.Range("N62").Value = -333341918
$ws.Range("L62").Value = 333340670
$ws.Range("M65").Value = -30321.11
$ws.Range("N65").Value = -1666709590
$ws.Range("H65").Value = 83340184
$ws.Range("I65").Value = 6688.222
$ws.Range("L65").Value = 1666703350
$ws.Range("J65").Value = 333340670
$ws.Range("K65").Value = 33441.11
$ws.Range("M99").Value = -431.8
$ws.Range("I99").Value = 1929.8
$ws.Range("H99").Value = 2395.5908
$ws.Range("J99").Value = 3393.7144
$ws.Range("N99").Value = -6389.7144
$ws.Range("L99").Value = 3393.7144
$ws.Range("K99").Value = 1929.8
$ws.Range("I126").Value = 1929.8
$ws.Range("H126").Value = 2395.5908
$ws.Range("J126").Value = 3393.7144
$ws.Range("L126").Value = 10181.1432
$ws.Range("K126").Value = 5789.4
$ws.Range("N126").Value = -15121.1432
$ws.Range("M126").Value = -3319.4
$ws.Range("K134").Value = 4846.5
$ws.Range("H134").Value = 3067.8235
$ws.Range("I134").Value = 1615.5
$ws.Range("M134").Value = -2311.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("L17").Value = 2441.25
$ws.Range("J17").Value = 813.75
$ws.Range("I17").Value = 609.75
$ws.Range("K17").Value = 1829.25
$ws.Range("N17").Value = -2779.25
$ws.Range("M17").Value = -1660.25
$ws.Range("H17").Value = 745.75
$ws.Range("L122").Value = 6678
$ws.Range("N122").Value = -11578
$ws.Range("H122").Value = 612.8333
$ws.Range("J122").Value = 742

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("L80").Value = 4462.75
$ws.Range("I80").Value = 3990
$ws.Range("K80").Value = 3990
$ws.Range("H80").Value = 4368.2
$ws.Range("J80").Value = 4462.75
$ws.Range("N80").Value = -6458.75
$ws.Range("M80").Value = -2992
$ws.Range("J83").Value = 4462.75
$ws.Range("H83").Value = 4368.2
$ws.Range("M83").Value = -14958
$ws.Range("I83").Value = 3990
$ws.Range("L83").Value = 22313.75
$ws.Range("K83").Value = 19950
$ws.Range("N83").Value = -32297.75
$ws.Range("I97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("H97").Value = 0
$ws.Range("N97").ClearContents()
$ws.Range("K97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("M97").ClearContents()
$ws.Range("L102").Value = 4111.4
$ws.Range("I102").Value = 1289.1428
$ws.Range("N102").Value = -7355.4
$ws.Range("H102").Value = 2465.0833
$ws.Range("K102").Value = 1289.1428
$ws.Range("M102").Value = 332.8571999999999
$ws.Range("J102").Value = 4111.4
$ws.Range("I122").Value = 2376.5
$ws.Range("H122").Value = 71431360
$ws.Range("M122").Value = -4679.5
$ws.Range("K122").Value = 7129.5
$ws.Range("N132").Value = -20309.2505
$ws.Range("I132").Value = 2503.5334
$ws.Range("K132").Value = 7510.600199999999
$ws.Range("M132").Value = -4980.600199999999
$ws.Range("J132").Value = 5083.0835
$ws.Range("H132").Value = 3240.5476
$ws.Range("L132").Value = 15249.2505

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("M16").Value = -312
$ws.Range("H16").Value = 482
$ws.Range("I16").Value = 482
$ws.Range("K16").Value = 482
$ws.Range("L68").Value = 22975.5
$ws.Range("H68").Value = 10254.728
$ws.Range("J68").Value = 22975.5
$ws.Range("N68").Value = -24473.5
$ws.Range("L71").Value = 114877.5
$ws.Range("N71").Value = -122365.5
$ws.Range("J71").Value = 22975.5
$ws.Range("H71").Value = 10254.728
$ws.Range("I122").Value = 2932.1333
$ws.Range("H122").Value = 3713.4285
$ws.Range("M122").Value = -6346.3999
$ws.Range("K122").Value = 8796.3999
$ws.Range("I132").Value = 3772.7693
$ws.Range("K132").Value = 11318.3079
$ws.Range("M132").Value = -8788.3079
$ws.Range("H132").Value = 6357.44

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("K96").Value = 2273.75
$ws.Range("L96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("M96").Value = -900.75
$ws.Range("N96").ClearContents()
$ws.Range("I96").Value = 2273.75
$ws.Range("H96").Value = 2273.75
$ws.Range("K136").Value = 2889.6249
$ws.Range("I136").Value = 963.2083
$ws.Range("H136").Value = 2537
$ws.Range("M136").Value = -339.6248999999998
